$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated error-table values and N-counts (rows 2-10, columns B-G)
$data = @{
    2  = @(0.02300404818780313, 0.337879676387987, 0.2668330530619737, 0.5165588573066711, 0.5301872556910642, 19)
    3  = @(0.4362183256331637, 0.5571722338028687, 0.8814962611585107, 0.9388803231288376, 0.855493562520063, 18)
    4  = @(0.7297665574417316, 0.7297665574417316, 0.9114036334090323, 0.9546746217476572, 0.6344463573575141, 17)
    5  = @(0.6329849386531881, 0.6568932063127731, 0.6561129869449028, 0.8100080165929858, 0.5219890720522772, 16)
    6  = @(0.5673547796451167, 0.5928185759272669, 0.4654632102043935, 0.6822486425082819, 0.3922077131238274, 15)
    7  = @(0.5117157535595949, 0.5117157535595949, 0.3387636996991033, 0.5820341052714207, 0.287796619108349, 14)
    8  = @(0.4402895787565182, 0.4402895787565182, 0.2500964150522257, 0.5000964057581555, 0.2468365728875135, 13)
    9  = @(0.4177270777543842, 0.4177270777543842, 0.2295217317045278, 0.4790842636786641, 0.2488052405843257, 9)
    10 = @(0.2889056645135138, 0.2889056645135138, 0.1117927209483371, 0.3343541848823446, 0.188169597572051, 5)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
}
